# Add additional quality-grading columns (foreign matter / damage / weevil
# breakdown) to the "Paddy" and "Jowr" sheets, mirroring the richer column
# set already present on "Ragi".

$wb = $excel.ActiveWorkbook

# ---- Paddy sheet -----------------------------------------------------
$paddy = $wb.Worksheets.Item("Paddy")

# Split "Foreign Matter" into Organic / Inorganic, then append the new
# grading columns.
$paddy.Range("N1").Value = "Foreign Matter Organic"
$paddy.Range("O1").Value = "Foreign Matter Inorganic"
$paddy.Range("P1").Value = "Damaged Discoloured Sprouted And Weevilled Grains"
$paddy.Range("Q1").Value = "Immature Shrunken And Shriveled Grains"
$paddy.Range("R1").Value = "Admixture Of Lower Class"
$paddy.Range("S1").Value = "Amount"
$paddy.Range("T1").Value = "Remark"

$paddy.Range("P2").Value = "null"
$paddy.Range("Q2").Value = "null"
$paddy.Range("R2").Value = "null"
$paddy.Range("S2").Value = "null"
$paddy.Range("T2").Value = "null"

# ---- Jowr sheet --------------------------------------------------------
$jowr = $wb.Worksheets.Item("Jowr")

$jowr.Range("O1").Value = "Damaged Grains"
$jowr.Range("P1").Value = "Slightly Damaged"
$jowr.Range("Q1").Value = "Other Food Grains"
$jowr.Range("R1").Value = "Shriveled Immature Grains"
$jowr.Range("S1").Value = "Weevilled Grains"
$jowr.Range("T1").Value = "Amount"
$jowr.Range("U1").Value = "Remark"

$jowr.Range("O2").Value = "null"
$jowr.Range("P2").Value = "null"
$jowr.Range("Q2").Value = "null"
$jowr.Range("R2").Value = "null"
$jowr.Range("S2").Value = "null"
$jowr.Range("T2").Value = "null"
$jowr.Range("U2").Value = "null"
